# Flights workbook update:
#  - Remove the leftover empty placeholder row (row 13, cell H13) that had
#    no real data, just a cell style.
#  - Add the new flight record (RF55, United, San Salvador -> San Pedro
#    Sula) as row 11, right after the existing flight entries.
#  - Update the sheet's active selection to reflect the newly entered data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old, empty placeholder row 13 (only had a style on H13) ---
$ws.Rows.Item(13).Delete()

# --- Add the new flight entry in row 11 ---
$ws.Cells.Item(11, 1).Value() = 10
$ws.Cells.Item(11, 2).Value() = "United"
$ws.Cells.Item(11, 3).Value() = "RF55"
$ws.Cells.Item(11, 4).Value() = "San Salvador"
$ws.Cells.Item(11, 5).Value() = "San Pedro Sula"

# FECHA (date) column holds plain text like the other rows (e.g. "2/5/21"),
# so write it through a formula + paste-as-values round trip to avoid Excel
# auto-converting the typed text into a real date serial number.
$ws.Cells.Item(11, 6).Formula() = '="01/02/2021"'
$ws.Cells.Item(11, 6).Copy()
$ws.Cells.Item(11, 6).PasteSpecial(-4163)
$excel.CutCopyMode() = 0

$ws.Cells.Item(11, 7).Value() = "19:30"
$ws.Cells.Item(11, 8).Value() = "21:15"
$ws.Cells.Item(11, 9).Value() = "On time"
$ws.Cells.Item(11, 10).Value() = "No comments yet"

# --- Update the active selection/view to the newly entered row ---
$ws.Activate()
$null = $ws.Range("A11:U230").Select()
